$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.051931324470075
$ws.Range("D2").Value = 0.0334777856361157
$ws.Range("E2").Value = 0.252259609794887
$ws.Range("F2").Value = 1.035586324422141
$ws.Range("G2").Value = 0.8951436972055831
$ws.Range("H2").Value = 0.918194064002293
$ws.Range("K2").Value = 0.4958702926694514
$ws.Range("L2").Value = 0.1635153336092543
$ws.Range("N2").Value = 2.053928379703617
$ws.Range("B3").Value = 1.025680380202431
$ws.Range("D3").Value = 0.03298472453969836
$ws.Range("E3").Value = 0.2536670297256984
$ws.Range("F3").Value = 1.023944163161431
$ws.Range("G3").Value = 0.8842639870132984
$ws.Range("H3").Value = 0.9179156127566444
$ws.Range("K3").Value = 0.4330694660474421
$ws.Range("L3").Value = 0.1521611400681735
$ws.Range("N3").Value = 2.07398026373059
$ws.Range("B4").Value = 1.010096680189065
$ws.Range("D4").Value = 0.03267793606853786
$ws.Range("E4").Value = 0.2545929883362419
$ws.Range("F4").Value = 1.017430457176999
$ws.Range("G4").Value = 0.8781722781685204
$ws.Range("H4").Value = 0.9181918428848803
$ws.Range("K4").Value = 0.3944948579348875
$ws.Range("L4").Value = 0.1452711664136359
$ws.Range("N4").Value = 2.086956082072284
$ws.Range("B5").Value = 1.003880971223083
$ws.Range("D5").Value = 0.0325519082637058
$ws.Range("E5").Value = 0.2549858918794676
$ws.Range("F5").Value = 1.014935476847882
$ws.Range("G5").Value = 0.8758375533005278
$ws.Range("H5").Value = 0.91841684397356
$ws.Range("K5").Value = 0.3787718904165729
$ws.Range("L5").Value = 0.1424839735514354
$ws.Range("N5").Value = 2.092410688369132
$ws.Range("B6").Value = 1.002857008652427
$ws.Range("D6").Value = 0.03253092075684449
$ws.Range("E6").Value = 0.2550520743153202
$ws.Range("F6").Value = 1.014530810039595
$ws.Range("G6").Value = 0.8754587862923557
$ws.Range("H6").Value = 0.9184609955703706
$ws.Range("K6").Value = 0.3761608976111575
$ws.Range("L6").Value = 0.1420224027625068
$ws.Range("N6").Value = 2.093326496829476
$ws.Range("B7").Value = 1.010012306766839
$ws.Range("D7").Value = 0.03267624048763551
$ws.Range("E7").Value = 0.254598224098415
$ws.Range("F7").Value = 1.017396163797756
$ws.Range("G7").Value = 0.878140193611884
$ws.Range("H7").Value = 0.9181944221085416
$ws.Range("K7").Value = 0.3942828264299578
$ws.Range("L7").Value = 0.1452334942118227
$ws.Range("N7").Value = 2.087028969325956
$ws.Range("B8").Value = 1.042769279839177
$ws.Range("D8").Value = 0.03330862312976279
$ws.Range("E8").Value = 0.2527320793589354
$ws.Range("F8").Value = 1.031440278290134
$ws.Range("G8").Value = 0.8912700592390905
$ws.Range("H8").Value = 0.9180052245568646
$ws.Range("K8").Value = 0.4742196971004091
$ws.Range("L8").Value = 0.1595834916662398
$ws.Range("N8").Value = 2.060704432618181
$ws.Range("B9").Value = 1.111234607197048
$ws.Range("D9").Value = 0.03451628862999101
$ws.Range("E9").Value = 0.2495616301935977
$ws.Range("F9").Value = 1.064027044693617
$ws.Range("G9").Value = 0.9217023190303877
$ws.Range("H9").Value = 0.9211845049246676
$ws.Range("K9").Value = 0.6308593063733099
$ws.Range("L9").Value = 0.1883715054239161
$ws.Range("N9").Value = 2.014352682600379
$ws.Range("B10").Value = 1.164104922945256
$ws.Range("D10").Value = 0.0353834327200957
$ws.Range("E10").Value = 0.247528688935164
$ws.Range("F10").Value = 1.091064099534805
$ws.Range("G10").Value = 0.9469418002459236
$ws.Range("H10").Value = 0.9256891175686803
$ws.Range("K10").Value = 0.745882912672073
$ws.Range("L10").Value = 0.2099202325470344
$ws.Range("N10").Value = 1.983514197157731
$ws.Range("B11").Value = 1.188713273261442
$ws.Range("D11").Value = 0.03577347984428414
$ws.Range("E11").Value = 0.2466678389313888
$ws.Range("F11").Value = 1.104040535292569
$ws.Range("G11").Value = 0.9590550853949082
$ws.Range("H11").Value = 0.9282104079331361
$ws.Range("K11").Value = 0.7982003626091512
$ws.Range("L11").Value = 0.2198106222319467
$ws.Range("N11").Value = 1.970184000869352
$ws.Range("B12").Value = 1.198111692939079
$ws.Range("D12").Value = 0.03592053746516655
$ws.Range("E12").Value = 0.2463510242929665
$ws.Range("F12").Value = 1.109052018296353
$ws.Range("G12").Value = 0.9637332908164922
$ws.Range("H12").Value = 0.9292331056931573
$ws.Range("K12").Value = 0.8180105537371105
$ws.Range("L12").Value = 0.2235684925635724
$ws.Range("N12").Value = 1.965236731115439
$ws.Range("B13").Value = 1.196084033266914
$ws.Range("D13").Value = 0.03588889479348012
$ws.Range("E13").Value = 0.24641884852117
$ws.Range("F13").Value = 1.107968362047259
$ws.Range("G13").Value = 0.9627216956750431
$ws.Range("H13").Value = 0.9290098269868565
$ws.Range("K13").Value = 0.8137441357704347
$ws.Range("L13").Value = 0.2227586082276076
$ws.Range("N13").Value = 1.966297737998822
$ws.Range("B14").Value = 1.189484890331954
$ws.Range("D14").Value = 0.03578559133602255
$ws.Range("E14").Value = 0.2466415907490465
$ws.Range("F14").Value = 1.10445087595933
$ws.Range("G14").Value = 0.9594381351687389
$ws.Range("H14").Value = 0.9282931839900073
$ws.Range("K14").Value = 0.7998301861725849
$ws.Range("L14").Value = 0.2201195322245866
$ws.Range("N14").Value = 1.969774969186155
$ws.Range("B15").Value = 1.185453102111609
$ws.Range("D15").Value = 0.03572223075307335
$ws.Range("E15").Value = 0.2467792203458945
$ws.Range("F15").Value = 1.10230902947599
$ws.Range("G15").Value = 0.9574387419243067
$ws.Range("H15").Value = 0.9278630687232123
$ws.Range("K15").Value = 0.7913073151792673
$ws.Range("L15").Value = 0.218504661358395
$ws.Range("N15").Value = 1.971917977778958
$ws.Range("B16").Value = 1.162507892850925
$ws.Range("D16").Value = 0.03535785248633161
$ws.Range("E16").Value = 0.2475862321355482
$ws.Range("F16").Value = 1.090229703250813
$ws.Range("G16").Value = 0.9461629082926066
$ws.Range("H16").Value = 0.9255338497789865
$ws.Range("K16").Value = 0.7424636761648458
$ws.Range("L16").Value = 0.2092756361515029
$ws.Range("N16").Value = 1.984399423674432
$ws.Range("B17").Value = 1.14857423830594
$ws.Range("D17").Value = 0.03513317942672955
$ws.Range("E17").Value = 0.248097668016193
$ws.Range("F17").Value = 1.082993014105199
$ws.Range("G17").Value = 0.9394075802283766
$ws.Range("H17").Value = 0.9242259054644677
$ws.Range("K17").Value = 0.7124976334238511
$ws.Range("L17").Value = 0.2036363883576939
$ws.Range("N17").Value = 1.992235371992411
$ws.Range("B18").Value = 1.140612454433438
$ws.Range("D18").Value = 0.03500353782250798
$ws.Range("E18").Value = 0.2483978526952146
$ws.Range("F18").Value = 1.078894384941222
$ws.Range("G18").Value = 0.9355815349473318
$ws.Range("H18").Value = 0.9235180455563494
$ws.Range("K18").Value = 0.6952612787914063
$ws.Range("L18").Value = 0.2004011046634417
$ws.Range("N18").Value = 1.996808127788757
$ws.Range("B19").Value = 1.137925757245085
$ws.Range("D19").Value = 0.03495957229968383
$ws.Range("E19").Value = 0.2485005248137213
$ws.Range("F19").Value = 1.077517597457614
$ws.Range("G19").Value = 0.9342963008169534
$ws.Range("H19").Value = 0.9232860066321109
$ws.Range("K19").Value = 0.6894252356821085
$ws.Range("L19").Value = 0.1993071130812041
$ws.Range("N19").Value = 1.998367667972325
$ws.Range("B20").Value = 1.150052070293924
$ws.Range("D20").Value = 0.03515713933605014
$ws.Range("E20").Value = 0.248042601896227
$ws.Range("F20").Value = 1.083756775359234
$ws.Range("G20").Value = 0.9401205421705612
$ws.Range("H20").Value = 0.9243605391474858
$ws.Range("K20").Value = 0.7156876418382012
$ws.Range("L20").Value = 0.204235840967371
$ws.Range("N20").Value = 1.991394418416292
$ws.Range("B21").Value = 1.191421056617259
$ws.Range("D21").Value = 0.03581595163547746
$ws.Range("E21").Value = 0.2465759172662878
$ws.Range("F21").Value = 1.10548139664273
$ws.Range("G21").Value = 0.9604001196233298
$ws.Range("H21").Value = 0.9285018351415033
$ws.Range("K21").Value = 0.8039170873547334
$ws.Range("L21").Value = 0.2208943512190018
$ws.Range("N21").Value = 1.968750890927488
$ws.Range("B22").Value = 1.218922874744095
$ws.Range("D22").Value = 0.03624276130716808
$ws.Range("E22").Value = 0.2456707958864595
$ws.Range("F22").Value = 1.120248613127643
$ws.Range("G22").Value = 0.9741855403790396
$ws.Range("H22").Value = 0.9316044405861703
$ws.Range("K22").Value = 0.8615724939363929
$ws.Range("L22").Value = 0.2318550483911821
$ws.Range("N22").Value = 1.954538406886762
$ws.Range("B23").Value = 1.204202233072351
$ws.Range("D23").Value = 0.03601531214644638
$ws.Range("E23").Value = 0.2461489943731099
$ws.Range("F23").Value = 1.112314943977438
$ws.Range("G23").Value = 0.9667792573385725
$ws.Range("H23").Value = 0.9299122682893994
$ws.Range("K23").Value = 0.8308014937032056
$ws.Range("L23").Value = 0.2259984107720499
$ws.Range("N23").Value = 1.962070170680661
$ws.Range("B24").Value = 1.149383790230189
$ws.Range("D24").Value = 0.03514630853612033
$ws.Range("E24").Value = 0.2480674781211913
$ws.Range("F24").Value = 1.083411286281461
$ws.Range("G24").Value = 0.9397980324196311
$ws.Range("H24").Value = 0.9242995339045308
$ws.Range("K24").Value = 0.7142454654296273
$ws.Range("L24").Value = 0.2039648072772025
$ws.Range("N24").Value = 1.991774402445142
$ws.Range("B25").Value = 1.092261141298252
$ws.Range("D25").Value = 0.03419309261238723
$ws.Range("E25").Value = 0.2503671459129047
$ws.Range("F25").Value = 1.054669238504843
$ws.Range("G25").Value = 0.9129654809214003
$ws.Range("H25").Value = 0.9199438874669852
$ws.Range("K25").Value = 0.5884954219815199
$ws.Range("L25").Value = 0.1805139010115937
$ws.Range("N25").Value = 2.026327504213608
